$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value looks like a plain number (e.g. "1.004") must be
# pre-formatted as Text so Excel stores them as strings (matching the
# source workbook, which keeps these as text, not numbers). The style is
# reset back to Normal afterwards so no lasting style/number-format change
# is left on the cell.
$numericLookingCells = @(
    "D4",
    "D5",
    "D7",
    "D9",
    "D10",
    "D11",
    "D12",
    "D14",
    "D15",
    "D16",
    "D17",
    "D18",
    "D19",
    "D20",
    "D21",
    "D22",
    "D24",
    "D25",
    "D26",
    "D28",
    "D29",
    "D30",
    "D31",
    "D33",
    "D34",
    "D35",
    "D36",
    "D37",
    "D38",
    "D39",
    "D40",
    "D41",
    "D42",
    "D43",
    "D44",
    "D45",
    "D46",
    "D47",
    "D48",
    "D49",
    "D50",
    "D51",
)
foreach ($addr in $numericLookingCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated values
$ws.Range("D2").Value = "27.955.54"
$ws.Range("E2").Value = "  +0.01%  "
$ws.Range("D3").Value = "1.854.07"
$ws.Range("E3").Value = "  -0.70%  "
$ws.Range("D4").Value = "1.004"
$ws.Range("E4").Value = "  +0.22%  "
$ws.Range("D5").Value = "311.72"
$ws.Range("E6").Value = "  +0.18%  "
$ws.Range("D7").Value = "0.5063"
$ws.Range("E7").Value = "  +1.58%  "
$ws.Range("E8").Value = "  -0.28%  "
$ws.Range("D9").Value = "0.08225"
$ws.Range("E9").Value = "  -8.06%  "
$ws.Range("D10").Value = "41.55"
$ws.Range("E10").Value = "  +0.26%  "
$ws.Range("D11").Value = "1.105"
$ws.Range("E11").Value = "  -1.27%  "
$ws.Range("D12").Value = "6.187"
$ws.Range("D13").Value = "1.858.78"
$ws.Range("E13").Value = "  -0.60%  "
$ws.Range("D14").Value = "20.42"
$ws.Range("E14").Value = "  -1.24%  "
$ws.Range("D15").Value = "7.169"
$ws.Range("E15").Value = "  -0.91%  "
$ws.Range("D16").Value = "1.004"
$ws.Range("E16").Value = "  +0.22%  "
$ws.Range("D17").Value = "0.00001090"
$ws.Range("E17").Value = "  -0.91%  "
$ws.Range("D18").Value = "90.20"
$ws.Range("E18").Value = "  -0.65%  "
$ws.Range("D19").Value = "0.06601"
$ws.Range("E19").Value = "  -0.76%  "
$ws.Range("D20").Value = "17.65"
$ws.Range("E20").Value = "  -1.79%  "
$ws.Range("D21").Value = "1.003"
$ws.Range("E21").Value = "  +0.20%  "
$ws.Range("D22").Value = "5.992"
$ws.Range("E22").Value = "  -2.01%  "
$ws.Range("D23").Value = "27.978.39"
$ws.Range("D24").Value = "11.02"
$ws.Range("E24").Value = "  -4.87%  "
$ws.Range("D25").Value = "2.238"
$ws.Range("E25").Value = "  -2.14%  "
$ws.Range("D26").Value = "2.539"
$ws.Range("E26").Value = "  +1.23%  "
$ws.Range("D27").Value = "2.071.99"
$ws.Range("E27").Value = "  -0.66%  "
$ws.Range("D28").Value = "157.72"
$ws.Range("E28").Value = "  -0.47%  "
$ws.Range("D29").Value = "20.32"
$ws.Range("E29").Value = "  -1.74%  "
$ws.Range("D30").Value = "123.98"
$ws.Range("E30").Value = "  -1.84%  "
$ws.Range("D31").Value = "0.1055"
$ws.Range("E31").Value = "  -0.31%  "
$ws.Range("E32").Value = "  -2.29%  "
$ws.Range("D33").Value = "5.579"
$ws.Range("E33").Value = "  +0.09%  "
$ws.Range("D34").Value = "3.592"
$ws.Range("E34").Value = "  +0.34%  "
$ws.Range("D35").Value = "9.497"
$ws.Range("E35").Value = "  +0.77%  "
$ws.Range("D36").Value = "0.06518"
$ws.Range("E36").Value = "  -0.37%  "
$ws.Range("D37").Value = "0.02395"
$ws.Range("E37").Value = "  -0.24%  "
$ws.Range("D38").Value = "0.2156"
$ws.Range("E38").Value = "  -1.60%  "
$ws.Range("D39").Value = "1.193"
$ws.Range("E39").Value = "  -0.50%  "
$ws.Range("D40").Value = "1.235"
$ws.Range("E40").Value = "  -4.72%  "
$ws.Range("D41").Value = "0.6361"
$ws.Range("E41").Value = "  -0.11%  "
$ws.Range("B42").Value = "InternetComputer(DFINITY)"
$ws.Range("C42").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D42").Value = "4.844"
$ws.Range("E42").Value = "  -1.15%  "
$ws.Range("B43").Value = "Aptos"
$ws.Range("C43").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D43").Value = "11.12"
$ws.Range("E43").Value = "  -4.79%  "
$ws.Range("D44").Value = "0.6015"
$ws.Range("E44").Value = "  +0.22%  "
$ws.Range("D45").Value = "13.01"
$ws.Range("E45").Value = "  -1.15%  "
$ws.Range("D46").Value = "1.277"
$ws.Range("E46").Value = "  -0.06%  "
$ws.Range("D47").Value = "3.651"
$ws.Range("E47").Value = "  -0.63%  "
$ws.Range("D48").Value = "1.975"
$ws.Range("E48").Value = "  -0.46%  "
$ws.Range("D49").Value = "1.198"
$ws.Range("E49").Value = "  -2.13%  "
$ws.Range("D50").Value = "120.28"
$ws.Range("D51").Value = "78.32"
$ws.Range("E51").Value = "  +0.20%  "

# Restore default styling on the cells we text-formatted above
foreach ($addr in $numericLookingCells) {
    $ws.Range($addr).Style = "Normal"
}
